$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the extraction timestamp (column B, rows 2-8)
$newTimestamp = "26/01/2026 22:57:21"
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 2).Value = $newTimestamp
}

# 2. Row 2's "Documento" cell (N2) gets a document link for the first time.
#    Add the hyperlink first, then copy the hyperlink-style formatting already
#    used by N3 (blue/underline) onto N2 (Hyperlinks.Add applies its own
#    built-in style, so re-apply the sheet's existing link look afterwards)
#    and set the display text.
$ws.Hyperlinks.Add($ws.Range("N2"), "/home/sebastianrojas/Escritorio/Poc_prueba/documentos_descargados/CP-ABR-1-2026-MDT-1.pdf")
$ws.Range("N3").Copy()
$ws.Range("N2").PasteSpecial(-4122)
$ws.Range("N2").Value = "Ver Documento"

# 3. Rows 3-8 already had "Ver .PDF" links; just relabel the link text.
for ($r = 3; $r -le 8; $r++) {
    $ws.Cells.Item($r, 14).Value = "Ver Documento"
}
